$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "End Time" for the last existing record (row 127).
$ws.Range("C127").Value = 0.80555555555555547

# Add a new row for the table, pushing in an empty formatted row so the
# new cells inherit formatting (number formats, styles) from row 127,
# the same way Excel does when a table grows by one record.
$ws.Rows.Item(128).Insert()

# New daily record.
$ws.Range("A128").Value = 43450
$ws.Range("B128").Clear()
$ws.Range("C128").Clear()
$ws.Range("D128").Formula = "=(C128-B128)* 1440"
$ws.Range("E128").Formula = "=IF(C128>B128, (C128-B128)*1440, (B128-C128)*1440)"
$ws.Range("F128").Formula = "=ABS((C128-B128)*1440)"

# Grow the table (ListObject) so its range/autofilter covers the new row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F128"))

# Leave the selection where the user would be after entering the date,
# ready to type the next (Start Time) value.
$ws.Range("B128").Select() | Out-Null
